$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting for the new rows (43-49) by copying the
# existing date-formatted row so the new cells share style "s=2"
# instead of picking up a default/general format.
$ws.Range("A41:B41").Copy() | Out-Null
$ws.Range("A43:B49").PasteSpecial(-4122) | Out-Null

$data = @(
    @(2, 45207, 45211),
    @(3, 45220, 45222),
    @(4, 45227, 45231),
    @(5, 45234, 45235),
    @(6, 45238, 45239),
    @(7, 45243, 45247),
    @(8, 45255, 45256),
    @(9, 45264, 45270),
    @(10, 45271, 45273),
    @(11, 45276, 45278),
    @(12, 45281, 45284),
    @(13, 45289, 45290),
    @(14, 45297, 45299),
    @(15, 45302, 45303),
    @(16, 45308, 45312),
    @(17, 45319, 45320),
    @(18, 45334, 45336),
    @(19, 45339, 45344),
    @(20, 45349, 45354),
    @(21, 45355, 45361),
    @(22, 45362, 45364),
    @(23, 45367, 45371),
    @(24, 45376, 45382),
    @(25, 45383, 45389),
    @(26, 45394, 45398),
    @(27, 45402, 45403),
    @(28, 45404, 45410),
    @(29, 45411, 45417),
    @(30, 45418, 45424),
    @(31, 45425, 45428),
    @(32, 45440, 45445),
    @(33, 45446, 45452),
    @(34, 45455, 45459),
    @(35, 45460, 45466),
    @(36, 45467, 45473),
    @(37, 45474, 45480),
    @(38, 45485, 45487),
    @(39, 45488, 45494),
    @(40, 45495, 45501),
    @(41, 45502, 45508),
    @(42, 45509, 45515),
    @(43, 45516, 45522),
    @(44, 45523, 45529),
    @(45, 45530, 45536),
    @(46, 45537, 45543),
    @(47, 45544, 45550),
    @(48, 45551, 45557),
    @(49, 45558, 45563)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

Write-Output "updated rows"